# update code / notes
#
# Applies 5 text edits across 3 slides of the "problem-solving" deck:
#   Slide 12: title "continued" suffix; "reuse patterns" bullet gains map/filter
#   Slide 14: "Write a function..." bullet reworded
#   Slide 5:  "bound the maximum size" bullet gains a parenthetical;
#             "Look for keywords" bullet (3 runs) consolidated + extended
#   Slide 8:  "by hand" clause split out and bolded

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 12: "Steps in "solve" phase" -> "... continued"
# ---------------------------------------------------------------------
$slide12 = $p.Slides.Item(12)
$title12 = $slide12.Shapes.Item(1).TextFrame.TextRange
$title12.Text = "Steps in " + [char]8220 + "solve" + [char]8221 + " phase continued"

$body12 = $slide12.Shapes.Item(2).TextFrame.TextRange
$reusePara = $body12.Paragraphs(4, 1)
$reuseChars = $body12.Characters($reusePara.Start, $reusePara.Length)
$reuseChars.Text = "Look for and reuse familiar programming patterns like vector sum, min, sort, map, filter, and find"

# ---------------------------------------------------------------------
# Slide 14: "Write a function definition ..." -> "Write a function ..."
# ---------------------------------------------------------------------
$slide14 = $p.Slides.Item(14)
$body14 = $slide14.Shapes.Item(2).TextFrame.TextRange
$writeFnPara = $body14.Paragraphs(1, 1)
$writeFnChars = $body14.Characters($writeFnPara.Start, $writeFnPara.Length)
$writeFnChars.Text = "Write a function that takes your input as parameter(s)"

# ---------------------------------------------------------------------
# Slide 5: "Can you bound the maximum size of the input?" -> add clause
# ---------------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
$body5 = $slide5.Shapes.Item(2).TextFrame.TextRange
$boundPara = $body5.Paragraphs(7, 1)
$boundChars = $body5.Characters($boundPara.Start, $boundPara.Length)
$boundChars.Text = "Can you bound the maximum size of the input (e.g., to fit in memory)?"

# Slide 5: "Look for keywords like ..." -> merge runs + append more keywords
$keywordsPara = $body5.Paragraphs(12, 1)
$keywordsChars = $body5.Characters($keywordsPara.Start, $keywordsPara.Length)
$keywordsChars.Text = "Look for keywords like min, max, average, median, sort, argmax, sum, find, search, collect, filter out, select, compute, etc..."

# ---------------------------------------------------------------------
# Slide 8: bold "by hand" within the "walk through" sentence
# ---------------------------------------------------------------------
$slide8 = $p.Slides.Item(8)
$body8 = $slide8.Shapes.Item(2).TextFrame.TextRange
$walkPara = $body8.Paragraphs(3, 1)
$walkText = $walkPara.Text
$byHandOffset = $walkText.IndexOf("by hand")
$byHandRange = $body8.Characters($walkPara.Start + $byHandOffset, 7)
$byHandRange.Font.Bold = $true
